$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-18 04:43:32"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-18 04:43:16"
$wsZhCn.Range("K2").Value = "2016-10-18 04:44:07"

# de-de sheet: Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-18 04:44:30"
